$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.964.56'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.172.97'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.57'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.618'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.94'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -5.07%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.566'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.72'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.12'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0926'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.00%  '
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.96'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.495.73'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.857'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.32'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.179.88'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '40.878.36'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0938'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.16'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.33'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.10'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.66'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +18.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.81'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +5.84%  '
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('E29').Value = '  -3.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.16'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.46'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.44'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.122'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.37'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +6.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0719'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.54'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.33'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +9.46%  '
$ws.Range('E39').Value = '  +3.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0294'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +8.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.20'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.99'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +17.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.61'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.91'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.199'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.78'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.53'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.100'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.01'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.12'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.85%  '
$ws.Range('E51').Value = '  -1.62%  '
